$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1541.9434
$ws.Range("J17").Value = 1266.6078
$ws.Range("L17").Value = 3799.8234
$ws.Range("N17").Value = -4135.8234
$ws.Range("H132").Value = 792.4375
$ws.Range("I132").Value = 738.76666
$ws.Range("J132").Value = 1597.5
$ws.Range("K132").Value = 2216.29998
$ws.Range("L132").Value = 4792.5
$ws.Range("M132").Value = 313.7000200000002
$ws.Range("N132").Value = -9852.5
$ws.Range("H137").Value = 1943.4117
$ws.Range("I137").Value = 1168.3334
$ws.Range("K137").Value = 3505.0002
$ws.Range("M137").Value = -955.0001999999999
$ws.Range("H140").Value = 55697.945
$ws.Range("J140").Value = 55697.945
$ws.Range("L140").Value = 55697.945
$ws.Range("N140").Value = -66057.94500000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3562.6736
$ws.Range("I32").Value = 3099.4783
$ws.Range("K32").Value = 3099.4783
$ws.Range("M32").Value = -2812.4783
$ws.Range("H61").Value = 2332.5557
$ws.Range("I61").Value = 1057.3334
$ws.Range("J61").Value = 6795.8335
$ws.Range("K61").Value = 1057.3334
$ws.Range("L61").Value = 6795.8335
$ws.Range("M61").Value = -845.3334
$ws.Range("N61").Value = -7219.8335
$ws.Range("H74").Value = 1684.1875
$ws.Range("I74").Value = 1852.2
$ws.Range("J74").Value = 1607.8182
$ws.Range("K74").Value = 1852.2
$ws.Range("L74").Value = 1607.8182
$ws.Range("M74").Value = -978.2
$ws.Range("N74").Value = -3355.8182
$ws.Range("H77").Value = 1684.1875
$ws.Range("I77").Value = 1852.2
$ws.Range("J77").Value = 1607.8182
$ws.Range("K77").Value = 9261
$ws.Range("L77").Value = 8039.090999999999
$ws.Range("M77").Value = -4893
$ws.Range("N77").Value = -16775.091
$ws.Range("H122").Value = 1182
$ws.Range("I122").Value = 1253.7646
$ws.Range("K122").Value = 3761.2938
$ws.Range("M122").Value = -1311.2938
$ws.Range("H132").Value = 1723.2927
$ws.Range("I132").Value = 1033.742
$ws.Range("J132").Value = 3860.9
$ws.Range("K132").Value = 3101.226
$ws.Range("L132").Value = 11582.7
$ws.Range("M132").Value = -571.2259999999997
$ws.Range("N132").Value = -16642.7
$ws.Range("H136").Value = 2332.5557
$ws.Range("I136").Value = 1057.3334
$ws.Range("J136").Value = 6795.8335
$ws.Range("K136").Value = 3172.0002
$ws.Range("L136").Value = 20387.5005
$ws.Range("M136").Value = -622.0001999999999
$ws.Range("N136").Value = -25487.5005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4219.17
$ws.Range("I134").Value = 4587.3257
$ws.Range("K134").Value = 13761.9771
$ws.Range("M134").Value = -11226.9771

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2300.2727
$ws.Range("I31").Value = 2860
$ws.Range("J31").Value = 2039.0667
$ws.Range("K31").Value = 2860
$ws.Range("L31").Value = 2039.0667
$ws.Range("M31").Value = -2565
$ws.Range("N31").Value = -2629.0667
$ws.Range("H34").Value = 2300.2727
$ws.Range("I34").Value = 2860
$ws.Range("J34").Value = 2039.0667
$ws.Range("K34").Value = 2860
$ws.Range("L34").Value = 2039.0667
$ws.Range("M34").Value = -2658
$ws.Range("N34").Value = -2443.0667

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 8487263
$ws.Range("J131").Value = 14053.302
$ws.Range("L131").Value = 42159.906
$ws.Range("N131").Value = -52239.906

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 8000
$ws.Range("I41").Value = 8000
$ws.Range("K41").Value = 8000
$ws.Range("M41").Value = -7645
$ws.Range("H70").Value = 4855.2856
$ws.Range("I70").Value = 5122
$ws.Range("J70").Value = 4499.6665
$ws.Range("K70").Value = 5122
$ws.Range("L70").Value = 4499.6665
$ws.Range("M70").Value = -4852
$ws.Range("N70").Value = -5039.6665
$ws.Range("H73").Value = 4855.2856
$ws.Range("I73").Value = 5122
$ws.Range("J73").Value = 4499.6665
$ws.Range("K73").Value = 5122
$ws.Range("L73").Value = 4499.6665
$ws.Range("M73").Value = -4186
$ws.Range("N73").Value = -6371.6665
$ws.Range("H132").Value = 940073.3
$ws.Range("I132").Value = 1166986.2
$ws.Range("J132").Value = 4057.125
$ws.Range("K132").Value = 3500958.6
$ws.Range("L132").Value = 12171.375
$ws.Range("M132").Value = -3498428.6
$ws.Range("N132").Value = -17231.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H124").Value = 14714
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 14714
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 14714
$ws.Range("N124").Value = -24534
$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("H127").Value = 75000
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 75000
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 75000
$ws.Range("N127").Value = -84920
$ws.Range("H128").Value = 44999
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 44999
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 44999
$ws.Range("N128").Value = -54959
$ws.Range("H129").Value = 0
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("H130").Value = 20000
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 20000
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 20000
$ws.Range("N130").Value = -30040
$ws.Range("H131").Value = 46130.4
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 46130.4
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 46130.4
$ws.Range("N131").Value = -56210.4
$ws.Range("H132").Value = 1836.7097
$ws.Range("I132").Value = 1122.5416
$ws.Range("J132").Value = 4285.2856
$ws.Range("K132").Value = 3367.6248
$ws.Range("L132").Value = 12855.8568
$ws.Range("M132").Value = -837.6248000000001
$ws.Range("N132").Value = -17915.8568
$ws.Range("H133").Value = 70326
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 70326
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 70326
$ws.Range("N133").Value = -75386
$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("H135").Value = 33557.25
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 33557.25
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 33557.25
$ws.Range("N135").Value = -43697.25
$ws.Range("H136").Value = 2133.3684
$ws.Range("I136").Value = 1202.8064
$ws.Range("J136").Value = 6254.4287
$ws.Range("K136").Value = 3608.4192
$ws.Range("L136").Value = 18763.2861
$ws.Range("M136").Value = -1058.4192
$ws.Range("N136").Value = -23863.2861
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("H138").Value = 0
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("H139").Value = 45000
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 45000
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 45000
$ws.Range("N139").Value = -55280
$ws.Range("H140").Value = 0
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("H141").Value = 48713
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 48713
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 48713
$ws.Range("N141").Value = -59073
